$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newTextBoth = "Automatic Connection Feature Completed (Both inputs and outputs)"
$newTextPorts = "Automatic Connection Feature Completed, Updated port names"

# Row 37 gets the "Updated port names" variant.
$ws.Cells.Item(37, 4).Value = $newTextPorts

# Rows where column D currently reads "Automatic Connection Feature Completed"
# (all except row 37) get the "(Both inputs and outputs)" suffix.
$rows = @(3,4,6,8,9,11,19,20,21,22,23,24,25,27,28,31,34,35)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 4).Value = $newTextBoth
}

# Selection moves from J37 to D37.
$ws.Range("D37").Select()

# Column D widens to fit the new, longer text (bestFit behaviour).
$ws.Columns.Item(4).ColumnWidth = 56.5
